$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '70.310.62'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.73%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.563.63'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.06%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '606.80'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +3.34%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '186.31'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.09%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '3.558.58'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.21%  '
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("E9").Value = '  -0.04%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.216'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +9.22%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.647'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.40%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '54.07'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("E13").Value = '  +1.74%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '9.48'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '4.128.95'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.01%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '70.307.09'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.76%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '3.578.97'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.66%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '12.71'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.65%  '
$ws.Range("E19").Value = '  -1.80%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '574.87'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +5.77%  '
$ws.Range("E21").Value = '  +0.78%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.995'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.66%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '17.38'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.64%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '4.73'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +3.67%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '4.96'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.44%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '93.92'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.77%  '
$ws.Range("E27").Value = '  -1.49%  '
$ws.Range("E28").Value = '  -1.87%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.36'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +3.03%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '32.38'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.82%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.06'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.83%  '
$ws.Range("E32").Value = '  -1.55%  '
$ws.Range("E33").Value = '  +2.13%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '64.04'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("E35").Value = '  +20.67%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '3.18'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.83%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '527.57'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.31%  '
$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.405'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.653.22'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +8.36%  '
$ws.Range("B40").Value = 'Dai'
$ws.Range("C40").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '37.54'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("E42").Value = '  +2.92%  '
$ws.Range("E43").Value = '  +4.20%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.138'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.48%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0457'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +3.99%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.48'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.87%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.95'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("E49").Value = '  +1.50%  '
$ws.Range("E50").Value = '  +0.35%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '135.64'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.90%  '
